$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Character used for the PEPE subscript-3 price notation (0.0₃...)
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.329.57'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -4.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.902.96'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.94'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.509'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.902.36'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("E10").Value = '  -6.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.72'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000214'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.378.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.893.94'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.52'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +6.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.261.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '405.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.84'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.670'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.82'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.68'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.26'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.99'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '24.67'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0978'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.917'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.44'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("E36").Value = '  -9.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '48.03'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("B38").Value = 'Cosmos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.19'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.71%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = ('0.0{0}0640' -f $sub3)
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("E41").Value = '  +4.04%  '
$ws.Range("E42").Value = '  -5.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '365.63'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.608.45'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '120.60'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.95'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.30%  '
